{"js": "// The target paragraph originally holds three runs that together read:\n//   \"  12  \" + \"Distance_Range_4_point_corre_function_average\" + \"\"\n// The edit merges them into a single run (keeping the existing\n// \"  12  Distance_Range_4_point_corre_function_average\" text) and appends\n// the two new parameter entries \"13 Emin (for Intel MKL)\" and\n// \"14 Emax (for Intel MKL)\" onto the same line/run.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nconst marker = \"Distance_Range_4_point_corre_function_average\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the paragraph containing '\" + marker + \"'\");\n}\n\n// Replacing the whole paragraph's text collapses the existing runs into a\n// single run carrying the paragraph's (shared) run formatting \u2014 matching\n// the diff, which leaves one run with the combined text and the original\n// JetBrains Mono / black rPr.\ntarget.insertText(\n  \"  12  Distance_Range_4_point_corre_function_average  13 Emin  (for Intel MKL)   14 Emax (for Intel MKL)\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# The target paragraph originally holds three runs whose text together\n# reads \"  12  \" + \"Distance_Range_4_point_corre_function_average\" + \"\".\n# The edit merges them into one run (keeping the original text) and\n# appends the two new parameter entries \"13 Emin (for Intel MKL)\" and\n# \"14 Emax (for Intel MKL)\" onto the same line/run.\n\n$d = $word.ActiveDocument\n\n$marker = \"Distance_Range_4_point_corre_function_average\"\n\n$r = $d.Content\n$found = $r.Find.Execute($marker)\nif (-not $found) {\n    throw \"Could not locate the paragraph containing '$marker'\"\n}\n\n# Expand the found range to its enclosing paragraph, then drop the\n# trailing paragraph mark so setting .Text doesn't clobber it.\n$r.Expand(4)  # wdParagraph\n$r.MoveEnd(1, -1)  # wdCharacter, shrink by one (the paragraph mark)\n\n$r.Text = \"  12  Distance_Range_4_point_corre_function_average  13 Emin  (for Intel MKL)   14 Emax (for Intel MKL)\"\n"}
